$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "merlin_bit"
$ws.Range("B6").Value = 720
$ws.Range("C6").Value = 3506
$ws.Range("D6").Value = 3357
$ws.Range("E6").Value = 0.1703738760056791
$ws.Range("F6").Value = 0.1766004415011038
$ws.Range("G6").Value = 0.1734312898952186
$ws.Range("H6").Value = 0.9305622528691291
